$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Block "Définition des tâches ... sprint 3" (rows 48-54) ---
# Row 48: new task entry (2h)
$ws.Range("B48").Value = "finalisaion des maquettes ainsi que l'ajout des maquette dans la documentation"
$ws.Range("D48").Value = 2

# Row 49: new task entry (1.5h)
$ws.Range("B49").Value = "Commencement du web summary"
$ws.Range("D49").Value = 1.5

# Row 54: personal reflection text for this half-day/day block
$ws.Range("B54").Value = "Je trouve que le groupe fonctionne très bien. La répartition à bien été effectuée. "

# Move the active selection to reflect where the author left off editing
$ws.Range("L64").Select()
